$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1868512110726644
$ws.Range("C2").Value = 0.5709342560553633
$ws.Range("J2").Value = 0.01730103806228374
$ws.Range("P2").Value = 0.1314878892733564
$ws.Range("S2").Value = 0.09342560553633218
$ws.Range("B3").Value = 0.0119047619047619
$ws.Range("C3").Value = 0.03571428571428571
$ws.Range("J3").Value = 0.01785714285714286
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.1845238095238095
$ws.Range("J4").Value = 0.06122448979591837
$ws.Range("P4").Value = 0.6530612244897959
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.0730593607305936
$ws.Range("D6").Value = 0.0136986301369863
$ws.Range("F6").Value = 0.0821917808219178
$ws.Range("J6").Value = 0.2237442922374429
$ws.Range("O6").Value = 0.0045662100456621
$ws.Range("Q6").Value = 0.1278538812785388
$ws.Range("R6").Value = 0.091324200913242
$ws.Range("S6").Value = 0.3835616438356164
$ws.Range("B7").Value = 0.1204188481675393
$ws.Range("D7").Value = 0.02094240837696335
$ws.Range("E7").Value = 0.005235602094240838
$ws.Range("F7").Value = 0.05759162303664921
$ws.Range("J7").Value = 0.1151832460732984
$ws.Range("O7").Value = 0.02094240837696335
$ws.Range("Q7").Value = 0.162303664921466
$ws.Range("R7").Value = 0.07853403141361257
$ws.Range("S7").Value = 0.418848167539267
$ws.Range("B8").Value = 0.09647058823529411
$ws.Range("D8").Value = 0.0188235294117647
$ws.Range("F8").Value = 0.04941176470588235
$ws.Range("J8").Value = 0.1341176470588235
$ws.Range("O8").Value = 0.01411764705882353
$ws.Range("Q8").Value = 0.1411764705882353
$ws.Range("R8").Value = 0.131764705882353
$ws.Range("S8").Value = 0.4141176470588235
$ws.Range("B9").Value = 0.1084745762711864
$ws.Range("D9").Value = 0.02372881355932203
$ws.Range("F9").Value = 0.06779661016949153
$ws.Range("J9").Value = 0.1220338983050848
$ws.Range("O9").Value = 0.01016949152542373
$ws.Range("Q9").Value = 0.1457627118644068
$ws.Range("R9").Value = 0.08813559322033898
$ws.Range("S9").Value = 0.4338983050847458
$ws.Range("B10").Value = 0.0893939393939394
$ws.Range("D10").Value = 0.02045454545454545
$ws.Range("E10").Value = 0.0007575757575757576
$ws.Range("F10").Value = 0.06363636363636363
$ws.Range("J10").Value = 0.1287878787878788
$ws.Range("O10").Value = 0.01893939393939394
$ws.Range("Q10").Value = 0.1795454545454545
$ws.Range("R10").Value = 0.1037878787878788
$ws.Range("S10").Value = 0.3946969696969697
$ws.Range("G11").Value = 0.1442953020134228
$ws.Range("J11").Value = 0.09731543624161074
$ws.Range("K11").Value = 0.2114093959731544
$ws.Range("L11").Value = 0.5335570469798657
$ws.Range("S11").Value = 0.01342281879194631
$ws.Range("G12").Value = 0.7804878048780488
$ws.Range("J12").Value = 0.1829268292682927
$ws.Range("L12").Value = 0.01219512195121951
$ws.Range("S12").Value = 0.02439024390243903
$ws.Range("G13").Value = 0.6585365853658537
$ws.Range("J13").Value = 0.2682926829268293
$ws.Range("S13").Value = 0.07317073170731707
$ws.Range("F15").Value = 0.02371541501976284
$ws.Range("H15").Value = 0.1422924901185771
$ws.Range("I15").Value = 0.06719367588932806
$ws.Range("J15").Value = 0.3438735177865613
$ws.Range("K15").Value = 0.03952569169960474
$ws.Range("M15").Value = 0.003952569169960474
$ws.Range("O15").Value = 0.1106719367588933
$ws.Range("S15").Value = 0.2687747035573123
$ws.Range("F16").Value = 0.005405405405405406
$ws.Range("H16").Value = 0.1783783783783784
$ws.Range("I16").Value = 0.08108108108108109
$ws.Range("J16").Value = 0.3837837837837838
$ws.Range("K16").Value = 0.1351351351351351
$ws.Range("M16").Value = 0.03243243243243243
$ws.Range("O16").Value = 0.07567567567567568
$ws.Range("S16").Value = 0.1081081081081081
$ws.Range("F17").Value = 0.02005012531328321
$ws.Range("H17").Value = 0.112781954887218
$ws.Range("I17").Value = 0.1629072681704261
$ws.Range("J17").Value = 0.4360902255639098
$ws.Range("K17").Value = 0.09774436090225563
$ws.Range("M17").Value = 0.01503759398496241
$ws.Range("O17").Value = 0.07017543859649122
$ws.Range("S17").Value = 0.08521303258145363
$ws.Range("F18").Value = 0.00784313725490196
$ws.Range("H18").Value = 0.1568627450980392
$ws.Range("I18").Value = 0.1411764705882353
$ws.Range("J18").Value = 0.3764705882352941
$ws.Range("K18").Value = 0.09803921568627451
$ws.Range("M18").Value = 0.01568627450980392
$ws.Range("O18").Value = 0.1098039215686274
$ws.Range("S18").Value = 0.09411764705882353
$ws.Range("F19").Value = 0.02119883040935672
$ws.Range("H19").Value = 0.2002923976608187
$ws.Range("I19").Value = 0.1184210526315789
$ws.Range("J19").Value = 0.3625730994152047
$ws.Range("K19").Value = 0.09868421052631579
$ws.Range("M19").Value = 0.01754385964912281
$ws.Range("N19").Value = 0.0007309941520467836
$ws.Range("O19").Value = 0.0577485380116959
$ws.Range("S19").Value = 0.1228070175438596
